$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (rename param, change value)
$ws.Range("A2").Value = "param_P_CHP1_max"
$ws.Range("B2").Value = 30

# Copy the formatting (style) of A2 down to the new A3:A5 cells before
# setting their values, so the new label cells pick up the same style
# (bordered, bold, centered) that A2 already has.
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A3").Value = "param_P_CHP1_min"
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = "param_P_heat_pump1_max"
$ws.Range("B4").Value = 50

$ws.Range("A5").Value = "param_P_heat_pump1_min"
$ws.Range("B5").Value = 0
